$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update D and E columns for simple value-only rows ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.907.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.51"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4824"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3823"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07383"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9400"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.11"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07807"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.898.20"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.492"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.597"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.74"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008874"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.011"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.953.73"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.81"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.126"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.134.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.81"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.939"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.24"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.052"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.74"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.982"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.330"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.221"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.730"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.135"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02049"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5654"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05391"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.991"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.068"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.559"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1532"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4904"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.69"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.011"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.87"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.675"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06107"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.81%  "

# --- Update E-only rows (no D change) ---
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.54%  "

# --- Rows 34 and 35 swap: Filecoin moves to row 34, ImmutableX moves to row 35 ---
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.659"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.95%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7625"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.31%  "
